# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "26.104.93"
$ws.Cells.Item(2,5).Value = "  +0.99%  "
$ws.Cells.Item(3,4).Value = "1.750.62"
$ws.Cells.Item(3,5).Value = "  +0.55%  "
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "1.000"
$ws.Cells.Item(4,5).Value = "  -0.06%  "
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "235.58"
$ws.Cells.Item(5,5).Value = "  +4.92%  "
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "1.001"
$ws.Cells.Item(6,5).Value = "  +0.01%  "
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "0.5306"
$ws.Cells.Item(7,5).Value = "  +3.00%  "
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.2797"
$ws.Cells.Item(8,5).Value = "  -0.06%  "
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "0.06194"
$ws.Cells.Item(9,5).Value = "  +1.87%  "
$ws.Cells.Item(10,4).Value = "1.744.78"
$ws.Cells.Item(10,5).Value = "  +0.30%  "
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.07182"
$ws.Cells.Item(11,5).Value = "  +3.28%  "
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "15.43"
$ws.Cells.Item(12,5).Value = "  +1.77%  "
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "0.6463"
$ws.Cells.Item(13,5).Value = "  +2.37%  "
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "4.623"
$ws.Cells.Item(14,5).Value = "  +3.12%  "
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "78.61"
$ws.Cells.Item(15,5).Value = "  +3.10%  "
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "1.001"
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "1.002"
$ws.Cells.Item(17,5).Value = "  +0.09%  "
$ws.Cells.Item(18,4).Value = "26.014.74"
$ws.Cells.Item(18,5).Value = "  +0.61%  "
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "11.70"
$ws.Cells.Item(19,5).Value = "  +2.72%  "
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "0.000006735"
$ws.Cells.Item(20,5).Value = "  +2.85%  "
$ws.Cells.Item(21,4).Value = "1.967.41"
$ws.Cells.Item(21,5).Value = "  +0.46%  "
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "4.329"
$ws.Cells.Item(22,5).Value = "  +6.38%  "
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "8.738"
$ws.Cells.Item(23,5).Value = "  +3.98%  "
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "5.249"
$ws.Cells.Item(24,5).Value = "  +2.97%  "
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "138.85"
$ws.Cells.Item(25,5).Value = "  +0.79%  "
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "1.516"
$ws.Cells.Item(26,5).Value = "  +0.46%  "
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "15.36"
$ws.Cells.Item(27,5).Value = "  +2.90%  "
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "1.807"
$ws.Cells.Item(28,5).Value = "  -0.43%  "
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "104.73"
$ws.Cells.Item(29,5).Value = "  +2.33%  "
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "0.08316"
$ws.Cells.Item(30,5).Value = "  +0.69%  "
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "3.807"
$ws.Cells.Item(31,5).Value = "  +5.60%  "
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "3.665"
$ws.Cells.Item(32,5).Value = "  +7.90%  "
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "0.04592"
$ws.Cells.Item(33,5).Value = "  +4.83%  "
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "2.648"
$ws.Cells.Item(34,5).Value = "  +0.95%  "
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = "1.011"
$ws.Cells.Item(35,5).Value = "  +5.01%  "
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "0.6364"
$ws.Cells.Item(36,5).Value = "  +6.27%  "
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = "2.719"
$ws.Cells.Item(37,5).Value = "  +1.93%  "
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "0.01610"
$ws.Cells.Item(38,5).Value = "  +4.09%  "
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "1.964"
$ws.Cells.Item(39,5).Value = "  +3.64%  "
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "1.001"
$ws.Cells.Item(40,5).Value = "  +0.04%  "
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "100.79"
$ws.Cells.Item(41,5).Value = "  +0.19%  "
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.3949"
$ws.Cells.Item(42,5).Value = "  +3.66%  "
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "0.7485"
$ws.Cells.Item(43,5).Value = "  +3.66%  "
$ws.Cells.Item(44,5).Value = "  +3.27%  "
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "0.1148"
$ws.Cells.Item(45,5).Value = "  +5.11%  "
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "6.370"
$ws.Cells.Item(46,5).Value = "  +1.85%  "
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "0.05357"
$ws.Cells.Item(47,5).Value = "  -1.80%  "
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "31.11"
$ws.Cells.Item(48,5).Value = "  +5.05%  "
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "54.54"
$ws.Cells.Item(49,5).Value = "  +4.77%  "

# Rows 50 and 51: coin identity swap (Decentraland <-> EnergySwap) plus independent Price/Volume updates
$ws.Cells.Item(50,2).Value = "EnergySwap"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "7.652"
# E50 stays "  +2.83%  " (unchanged)

$ws.Cells.Item(51,2).Value = "Decentraland"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "0.3476"
$ws.Cells.Item(51,5).Value = "  +3.35%  "

Write-Host "Update complete"
